# "more tests and improvements"
# Adds a new "agility" worksheet (placed right after "Sheet1") containing
# two small "benchmark" score tables (SF / LEUVEN columns), and makes the
# new sheet the active tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# New sheet goes right after Sheet1 and becomes the active tab.
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "agility"

# ---- Block 1: columns A:D ("SF" speed benchmark, top table) ----
$ws.Range("A1").Value = "BENCHMARK"
$ws.Range("B2").Value = "SPEED"
$ws.Range("A3").Value = "ACC"
$ws.Range("B3").Value = "LOW"
$ws.Range("C3").Value = "MED"
$ws.Range("D3").Value = "HIGH"
$ws.Range("A4").Value = "LOW"
$ws.Range("B4").Value = 33.16
$ws.Range("C4").Value = 18.33
$ws.Range("D4").Value = 132.07
$ws.Range("A5").Value = "MED"
$ws.Range("B5").Value = 8.91
$ws.Range("C5").Value = 17.2
$ws.Range("D5").Value = 34.7
$ws.Range("A6").Value = "HIGH"
$ws.Range("B6").Value = 7.35
$ws.Range("C6").Value = 5.65
$ws.Range("D6").Value = 4.82

# ---- Block 2: columns G:J ("SF" table) ----
$ws.Range("G1").Value = "SF"
$ws.Range("H2").Value = "SPEED"
$ws.Range("G3").Value = "ACC"
$ws.Range("H3").Value = "LOW"
$ws.Range("I3").Value = "MED"
$ws.Range("J3").Value = "HIGH"
$ws.Range("G4").Value = "LOW"
$ws.Range("H4").Value = 135.19
$ws.Range("I4").Value = 960.39
$ws.Range("J4").Value = "-"
$ws.Range("G5").Value = "MED"
$ws.Range("H5").Value = 37.74
$ws.Range("I5").Value = 65.84
$ws.Range("J5").Value = 800.51
$ws.Range("G6").Value = "HIGH"
$ws.Range("H6").Value = 25.9
$ws.Range("I6").Value = 22.24
$ws.Range("J6").Value = 93.68

# ---- Block 3: columns M:P ("LEUVEN" table) ----
$ws.Range("M1").Value = "LEUVEN"
$ws.Range("N2").Value = "SPEED"
$ws.Range("M3").Value = "ACC"
$ws.Range("N3").Value = "LOW"
$ws.Range("O3").Value = "MED"
$ws.Range("P3").Value = "HIGH"
$ws.Range("M4").Value = "LOW"
$ws.Range("N4").Value = 402.68
$ws.Range("O4").Value = "-"
$ws.Range("P4").Value = "-"
$ws.Range("M5").Value = "MED"
$ws.Range("N5").Value = 139.89
$ws.Range("O5").Value = 311.73
$ws.Range("P5").Value = "-"
$ws.Range("M6").Value = "HIGH"
$ws.Range("N6").Value = 99.5
$ws.Range("O6").Value = 119.17
$ws.Range("P6").Value = "-"

# ---- "SCORES" section header ----
$ws.Range("A9").Value = "SCORES"

# ---- Block 1b: columns A:D, second table ----
$ws.Range("A10").Value = "BENCHMARK"
$ws.Range("B11").Value = "SPEED"
$ws.Range("A12").Value = "ACC"
$ws.Range("B12").Value = "LOW"
$ws.Range("C12").Value = "MED"
$ws.Range("D12").Value = "HIGH"
$ws.Range("A13").Value = "LOW"
$ws.Range("B13").Value = 98.88
$ws.Range("C13").Value = 84.96
$ws.Range("D13").Value = 83.52
$ws.Range("A14").Value = "MED"
$ws.Range("B14").Value = 93.4
$ws.Range("C14").Value = 45.12
$ws.Range("D14").Value = 47.4
$ws.Range("A15").Value = "HIGH"
$ws.Range("B15").Value = 96
$ws.Range("C15").Value = 48.44
$ws.Range("D15").Value = 33.15

# ---- Block 2b: columns G:J, second table ----
$ws.Range("G10").Value = "SF"
$ws.Range("H11").Value = "SPEED"
$ws.Range("G12").Value = "ACC"
$ws.Range("H12").Value = "LOW"
$ws.Range("I12").Value = "MED"
$ws.Range("J12").Value = "HIGH"
$ws.Range("G13").Value = "LOW"
$ws.Range("H13").Value = 305.45
$ws.Range("I13").Value = 150.13
$ws.Range("J13").Value = "-"
$ws.Range("G14").Value = "MED"
$ws.Range("H14").Value = 317.32
$ws.Range("I14").Value = 107.72
$ws.Range("J14").Value = 74.2
$ws.Range("G15").Value = "HIGH"
$ws.Range("H15").Value = 321.1
$ws.Range("I15").Value = 108.52
$ws.Range("J15").Value = 61.8

# ---- Block 3b: columns M:P, second table ----
$ws.Range("M10").Value = "LEUVEN"
$ws.Range("N11").Value = "SPEED"
$ws.Range("M12").Value = "ACC"
$ws.Range("N12").Value = "LOW"
$ws.Range("O12").Value = "MED"
$ws.Range("P12").Value = "HIGH"
$ws.Range("M13").Value = "LOW"
$ws.Range("N13").Value = 280.92
$ws.Range("O13").Value = "-"
$ws.Range("P13").Value = "-"
$ws.Range("M14").Value = "MED"
$ws.Range("N14").Value = 284.4
$ws.Range("O14").Value = 99.6
$ws.Range("P14").Value = "-"
$ws.Range("M15").Value = "HIGH"
$ws.Range("N15").Value = 295.7
$ws.Range("O15").Value = 101.73
$ws.Range("P15").Value = "-"

# Match the author's final selection on the new sheet.
[void]$ws.Range("O16").Select()
